$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Carrera") values: "Negocios" (D2:D3) and "Manufactura" (D4:D5)
# both become "Sistemas".
$ws.Range("D2").Value = "Sistemas"
$ws.Range("D3").Value = "Sistemas"
$ws.Range("D4").Value = "Sistemas"
$ws.Range("D5").Value = "Sistemas"

# Update the active selection (cursor position) to reflect the saved view state.
$ws.Range("E8").Select()
